$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://24702 (94.16), maa://25390 (95.82), maa://36681 (86.3)"
$ws.Range("AB2").Value = "maa://21246 (91.32), maa://36684 (97.62), ***maa://22731 (6.67)"
$ws.Range("H3").Value = "maa://21247 (98.37), *maa://22748 (75.0)"
$ws.Range("L3").Value = "*maa://22880 (69.01), maa://20276 (84.11), *maa://22749 (66.67)"
$ws.Range("P3").Value = "maa://21249 (94.86), maa://26254 (95.83)"
$ws.Range("D5").Value = "maa://21245 (82.61), maa://22744 (83.33)"
$ws.Range("AF6").Value = "*maa://33152 (60.47), ***maa://22770 (27.27)"
$ws.Range("X7").Value = "maa://22399 (95.04), *maa://22758 (72.41)"
$ws.Range("AF7").Value = "*maa://26191 (68.42), *maa://36671 (70.21), *maa://42530 (63.64)"
$ws.Range("A8").Value = "更新日期：2024.11.19 13:18:32"
$ws.Range("P8").Value = "maa://32931 (84.85), *maa://21916 (60.66), maa://23252 (92.42), maa://37496 (96.15), **maa://22759 (45.45)"
$ws.Range("AB9").Value = "maa://28711 (88.3), ***maa://22740 (5.88), **maa://39938 (47.62), **maa://27377 (46.15), ***maa://25174 (20.0), maa://40166 (90.91)"
$ws.Range("T10").Value = "maa://27395 (95.78), maa://22755 (87.39), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("X10").Value = "maa://22301 (97.57), maa://22726 (100.0)"
$ws.Range("L11").Value = "maa://21287 (88.17)"
$ws.Range("T11").Value = "maa://22747 (93.2), maa://22501 (98.28)"
$ws.Range("X11").Value = "maa://36713 (98.08)"
$ws.Range("H12").Value = "maa://21867 (89.94)"
$ws.Range("X12").Value = "maa://22753 (91.61), *maa://21485 (76.87), maa://37962 (86.96)"
$ws.Range("D13").Value = "maa://24999 (91.59), maa://36673 (92.42), maa://25001 (85.51)"
$ws.Range("P13").Value = "maa://22676 (91.43), *maa://22583 (75.41), *maa://22500 (56.82)"
$ws.Range("X13").Value = "*maa://34957 (77.97), *maa://22768 (51.61)"
$ws.Range("AF13").Value = "**maa://22737 (30.37), maa://39883 (91.3), *maa://39885 (56.0)"
$ws.Range("L14").Value = "maa://26245 (96.24), maa://21288 (96.21), maa://36682 (97.3), maa://39841 (93.94)"
$ws.Range("AB14").Value = "maa://22764 (96.72)"
$ws.Range("D15").Value = "*maa://22743 (77.13), maa://22734 (83.76), *maa://30808 (63.93), ***maa://36048 (23.08)"
$ws.Range("H15").Value = "maa://24304 (88.72), maa://21478 (91.18)"
$ws.Range("D16").Value = "maa://21441 (96.24), maa://36679 (92.68), maa://37650 (96.77)"
$ws.Range("H17").Value = "maa://22430 (88.4), maa://39599 (84.38)"
$ws.Range("T17").Value = "**maa://42324 (40.0)"
$ws.Range("H18").Value = "maa://24421 (89.91)"
$ws.Range("T19").Value = "maa://24386 (98.94)"
$ws.Range("D20").Value = "maa://21432 (90.78), maa://25198 (92.86), *maa://20795 (50.4), maa://36680 (96.43)"
$ws.Range("H20").Value = "maa://22864 (88.57)"
$ws.Range("L20").Value = "maa://41331 (82.72)"
$ws.Range("D21").Value = "maa://21261 (97.37)"
$ws.Range("H22").Value = "maa://25236 (96.34), **maa://21678 (48.94), **maa://22735 (42.86)"
$ws.Range("X22").Value = "maa://21282 (98.41), *maa://37649 (69.57)"
$ws.Range("D23").Value = "***maa://28036 (28.36), **maa://41753 (50.0)"
$ws.Range("L23").Value = "maa://39756 (92.99), maa://39875 (93.22)"
$ws.Range("D24").Value = "maa://24368 (80.06)"
$ws.Range("X24").Value = "maa://29988 (86.36), maa://23504 (92.97), **maa://22892 (39.86), *maa://25141 (77.42), maa://36663 (80.95), ***maa://22815 (23.08)"
$ws.Range("D25").Value = "maa://29753 (95.08)"
$ws.Range("H25").Value = "*maa://29063 (73.79), *maa://25311 (74.49), ***maa://22725 (4.84)"
$ws.Range("L25").Value = "maa://24378 (86.84)"
$ws.Range("AB26").Value = "maa://42235 (91.8)"
$ws.Range("L28").Value = "*maa://30770 (79.55)"
$ws.Range("T28").Value = "maa://23263 (94.85), *maa://29765 (61.33)"
$ws.Range("X28").Value = "maa://39929 (89.2), ***maa://39723 (14.29), maa://41749 (85.29)"
$ws.Range("AF28").Value = "maa://36660 (92.47), *maa://36701 (62.96)"
$ws.Range("L29").Value = "maa://28432 (93.38), *maa://28440 (73.49), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("AF29").Value = "*maa://24080 (69.41), ***maa://34960 (8.7), maa://42865 (85.19)"
$ws.Range("AB30").Value = "maa://42979 (96.97)"
$ws.Range("L31").Value = "maa://35926 (93.82), maa://36258 (81.4)"
$ws.Range("L32").Value = "maa://28065 (94.87)"
$ws.Range("T32").Value = "maa://41108 (87.5), maa://42859 (93.62), maa://41238 (94.74)"
$ws.Range("AF32").Value = "maa://42408 (85.71)"
$ws.Range("P33").Value = "*maa://21956 (79.41), *maa://22730 (79.31)"
$ws.Range("L35").Value = "maa://41296 (95.7)"
$ws.Range("H36").Value = "maa://24375 (92.5)"
$ws.Range("P37").Value = "maa://21280 (89.18), *maa://21239 (72.73)"
$ws.Range("P38").Value = "*maa://24383 (68.42)"
$ws.Range("AF38").Value = "maa://36697 (85.53)"
$ws.Range("H39").Value = "maa://25199 (85.32), maa://36670 (88.16), maa://30434 (88.33), ***maa://25036 (16.0)"
$ws.Range("P39").Value = "maa://24709 (91.3)"
$ws.Range("P40").Value = "maa://23278 (95.93), maa://21386 (95.7), maa://36664 (90.2)"
$ws.Range("H41").Value = "maa://24466 (93.18)"
$ws.Range("H43").Value = "maa://22525 (92.31), maa://21284 (83.33)"
$ws.Range("H44").Value = "maa://29768 (97.75), maa://27728 (96.0)"
$ws.Range("H45").Value = "maa://21229 (85.08), maa://30807 (95.24), *maa://22767 (57.89), ***maa://20796 (13.79), *maa://42459 (66.67)"
$ws.Range("P45").Value = "*maa://36237 (61.54)"
$ws.Range("T45").Value = "**maa://39364 (38.89)"
$ws.Range("H46").Value = "maa://35931 (92.4)"
$ws.Range("H47").Value = "maa://27410 (96.0), maa://29661 (97.78), maa://28038 (84.62)"
$ws.Range("P49").Value = "*maa://39643 (68.18)"
$ws.Range("H51").Value = "maa://30769 (81.25)"
$ws.Range("H53").Value = "maa://32534 (93.29), **maa://32434 (34.78)"
$ws.Range("H55").Value = "maa://32532 (92.28)"
$ws.Range("H59").Value = "maa://27746 (82.69), maa://31270 (95.54)"
$ws.Range("H60").Value = "*maa://40438 (55.26)"
